# Correction type pour génération à partir fsh ea4a6f04ed193a83290686b2f69a3f9cd2e7f4ad
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-18T06:40:38+00:00"

# --- Elements sheet: fix the root Extension row's Short/Definition text
#     and fill in the (previously empty) RIM Mapping column ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("L2").Value = "Optional Extensions Element"
$elements.Range("M2").Value = "Optional Extension Element - found in all resources."
$elements.Range("AK2").Value = "N/A"
